# Updated capital structure database
# Applies refreshed financial metrics for the Guernsey Metals & Mining
# dataset (rows 2-5), including a company-name swap between rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = -2.285294117647059
$ws.Range("H2").Value = -2.285294117647059
$ws.Range("I2").Value = -2.649999999999999
$ws.Range("J2").Value = -2.649999999999999
$ws.Range("K2").Value = -9.27
$ws.Range("L2").Value = -2.726470588235294
$ws.Range("U2").Value = 3.974
$ws.Range("V2").Value = 0.03696744186046512
$ws.Range("W2").Value = -0.4235560588901472
$ws.Range("X2").Value = 0.04868667249302855
$ws.Range("Y2").Value = -0.4722427313831758
$ws.Range("Z2").Value = 0.1369918207824651
$ws.Range("AA2").Value = -0.1351610429447853
$ws.Range("AB2").Value = 0.04853708639830287
$ws.Range("AC2").Value = -0.1855235870821721
$ws.Range("AD2").Value = 3.91
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 3.91
$ws.Range("AG2").Value = -0.0640000000000005
$ws.Range("AH2").Value = 0.03509559285521946
$ws.Range("AI2").Value = 0.1100726310455492
$ws.Range("AJ2").Value = -0.0005957034885885597
$ws.Range("AK2").Value = -0.002028654748320036
$ws.Range("AL2").Value = 0.138
$ws.Range("AM2").Value = 0.138
$ws.Range("AN2").Value = -0.4671445639187574
$ws.Range("AO2").Value = -65.28985507246377
$ws.Range("AP2").Value = 0.007646356033452866
$ws.Range("AQ2").Value = -65.28985507246377

# Row 3
$ws.Range("K3").Value = -2.65
$ws.Range("U3").Value = 0.098
$ws.Range("V3").Value = 0.002824207492795389
$ws.Range("W3").Value = -4.26731078904992
$ws.Range("X3").Value = 0.04842924326104797
$ws.Range("Y3").Value = -4.315740032310968
$ws.Range("AA3").Value = 10.3921568627451
$ws.Range("AB3").Value = 0.04842924326104797
$ws.Range("AC3").Value = 10.34372761948405
$ws.Range("AG3").Value = -0.098
$ws.Range("AI3").Value = -0
$ws.Range("AJ3").Value = -0.002832206230853708
$ws.Range("AK3").Value = 0.1863117870722434
$ws.Range("AM3").Value = 0
$ws.Range("AP3").Value = 0.03712121212121212
$ws.Range("AQ3").ClearContents()

# Row 4
$ws.Range("B4").Value = "AfriTin Mining Limited (AIM:ATM)"
$ws.Range("G4").Value = -1.428571428571429
$ws.Range("H4").Value = -1.428571428571429
$ws.Range("I4").Value = -1.831168831168831
$ws.Range("J4").Value = -1.831168831168831
$ws.Range("K4").Value = -2.88
$ws.Range("L4").Value = -1.87012987012987
$ws.Range("O4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("U4").Value = 3.45
$ws.Range("V4").Value = 0.1268382352941177
$ws.Range("W4").Value = -0.1398058252427184
$ws.Range("X4").Value = 0.05362249944385618
$ws.Range("Y4").Value = -0.1934283246865746
$ws.Range("Z4").Value = 0.07381134969325154
$ws.Range("AA4").Value = -0.1351610429447853
$ws.Range("AB4").Value = 0.05036254413738683
$ws.Range("AC4").Value = -0.1855235870821721
$ws.Range("AD4").Value = 3.61
$ws.Range("AF4").Value = 3.61
$ws.Range("AG4").Value = 0.1599999999999997
$ws.Range("AH4").Value = 0.1171697500811425
$ws.Range("AI4").Value = 0.1187109503452811
$ws.Range("AJ4").Value = 0.005847953216374258
$ws.Range("AK4").Value = 0.005934718100890197
$ws.Range("AL4").Value = 0.128
$ws.Range("AM4").Value = 0.128
$ws.Range("AN4").Value = -1.388461538461538
$ws.Range("AO4").Value = -22.03125
$ws.Range("AP4").Value = -0.06153846153846142
$ws.Range("AQ4").Value = -22.03125

# Row 5
$ws.Range("B5").Value = "Ferro-Alloy Resources Limited (LSE:FAR)"
$ws.Range("G5").Value = -1.473118279569892
$ws.Range("H5").Value = -1.473118279569892
$ws.Range("I5").Value = -1.903225806451613
$ws.Range("J5").Value = -1.903225806451613
$ws.Range("K5").Value = -3.74
$ws.Range("L5").Value = -2.010752688172043
$ws.Range("U5").Value = 0.426
$ws.Range("V5").Value = 0.009342105263157895
$ws.Range("W5").Value = -0.4235560588901472
$ws.Range("X5").Value = 0.04868667249302855
$ws.Range("Y5").Value = -0.4722427313831758
$ws.Range("Z5").Value = 0.4418052256532067
$ws.Range("AA5").Value = -0.8408551068883611
$ws.Range("AB5").Value = 0.04853708639830287
$ws.Range("AC5").Value = -0.889392193286664
$ws.Range("AD5").Value = 0.3
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0.3
$ws.Range("AG5").Value = -0.126
$ws.Range("AH5").Value = 0.006535947712418301
$ws.Range("AI5").Value = 0.05415162454873646
$ws.Range("AJ5").Value = -0.002770814091568808
$ws.Range("AK5").Value = -0.02463824794681267
$ws.Range("AL5").Value = 0.01
$ws.Range("AM5").Value = 0.01
$ws.Range("AN5").Value = -0.09584664536741214
$ws.Range("AO5").Value = -354
$ws.Range("AP5").Value = 0.0402555910543131
$ws.Range("AQ5").Value = -354
